$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for 'biosat' (row 7) and 'O2_Ar_ratio' (row 8), which shifts
# the 'ncp' and 'k' rows up by two (former rows 9-10 become rows 7-8).
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()

# Update the active selection to match the post-edit state.
$ws.Range("A7:XFD8").Select()
